$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G (G2:G6): replace the per-row "out data" values (a formula date,
# plain numbers, an empty cell and a stray "aqui" string) with the single
# text value "22/04/2020", and give the column its medium-right outer border.
$ws.Range("G2:G6").NumberFormat = "General"
$ws.Range("G2:G6").Borders(10).LineStyle = 1
$ws.Range("G2:G6").Borders(10).Weight = -4138
$ws.Range("G2:G6").Value = "22/04/2020"

# --- Column A (A2:A6): give the scenario-id column its medium-left outer
# border all the way down, including the last row.
$ws.Range("A2:A6").Borders(7).LineStyle = 1
$ws.Range("A2:A6").Borders(7).Weight = -4138

# --- Row 6 (A6:G6): this is the last row of the table, so it additionally
# gets a medium bottom border across every column (on top of the existing
# left/right outer borders already applied to columns A and G above).
$ws.Range("A6:G6").Borders(9).LineStyle = 1
$ws.Range("A6:G6").Borders(9).Weight = -4138

# --- Selection: the saved workbook now has G2 selected instead of B5.
$ws.Range("G2").Select()
